$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: heading "An alternative way to work, using tabs"
#         -> "An Alternative Way to Work, Using Tabs" (title case)
# ---------------------------------------------------------------------------
$headingOld = "An alternative way to work, using tabs"
$headingNew = "An Alternative Way to Work, Using Tabs"

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*$headingOld*") {
        $p.Range.Find.Execute($headingOld, $true, $false, $false, $false, $false, $true, 1, $false, $headingNew, 2) | Out-Null
        break
    }
}

# ---------------------------------------------------------------------------
# Edit 2: paragraph about culture suffixes - append new sentence after
#         "...a single culture per language."
# ---------------------------------------------------------------------------
$tailOld = ", as is the case anyway. Other than that, the package will function just like in the case of using a single culture per language."
$tailNew = ", as is the case anyway. Other than that, the package will function just like in the case of using a single culture per language. Note however that, if tabs are used for translations, the suffix of each property has to include the culture, in capitals (e.g. bodyText_en-GB)."

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*$tailOld*") {
        $p.Range.Find.Execute($tailOld, $true, $false, $false, $false, $false, $true, 1, $false, $tailNew, 2) | Out-Null
        break
    }
}
